$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "ESFJ"
$ws.Range("B4").Value = "ENFJ"
$ws.Range("C4").Value = "ENFJ"

$ws.Range("A5").Value = "ISTJ"
$ws.Range("B5").Value = "INTP"
$ws.Range("C5").Value = "INTJ"

$ws.Range("A6").Value = "ISFP"
$ws.Range("B6").Value = "ISTP"
$ws.Range("C6").Value = "ISFP"

$ws.Range("E13").Select()
